$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is a blank separator row, like existing row 6 - copy its (empty) formatting/content
# down to row 19 so the row is materialised in the used range.
$ws.Range("A6:C6").Copy($ws.Range("A19:C19"))

# Row 20: SearchUser_SearchByDropdown
$ws.Cells.Item(20, 1).Value = "SearchUser_SearchByDropdown"
$ws.Cells.Item(20, 2).Value = "//*[@id=""root""]/div[1]/div/div[2]/div/div[2]/div/div[2]/div[1]/div[1]/div[1]/select"
$ws.Cells.Item(20, 3).Value = "By.xpath"

# Row 21: SearchUser_SearchBar
$ws.Cells.Item(21, 1).Value = "SearchUser_SearchBar"
$ws.Cells.Item(21, 2).Value = "//input[@placeholder='Search...']"
$ws.Cells.Item(21, 3).Value = "By.xpath"

# Row 22: SearchUser_SearchButton
$ws.Cells.Item(22, 1).Value = "SearchUser_SearchButton"
$ws.Cells.Item(22, 2).Value = "//*[@id=""root""]/div[1]/div/div[2]/div/div[2]/div/div[2]/div[1]/div[1]/div[2]/button"
$ws.Cells.Item(22, 3).Value = "By.xpath"

# Row 23: SearchUser_Result
$ws.Cells.Item(23, 1).Value = "SearchUser_Result"
$ws.Cells.Item(23, 2).Value = "//td[normalize-space()='Kasun Bandara']"
$ws.Cells.Item(23, 3).Value = "By.xpath"

# Rows 24 and 25 are trailing blank separator rows, same treatment as row 19.
$ws.Range("A6:C6").Copy($ws.Range("A24:C24"))
$ws.Range("A6:C6").Copy($ws.Range("A25:C25"))
